# Generate Report for Handoff
# Applies the localization-status.xlsx refresh: new source file names
# (ccbd7795... / ffff16413676...), new status text ("Ready for handoff"),
# refreshed handoff timestamps/xlf names, cleared "Latest Target
# File"/"Latest Handback File" columns, and the associated column-width
# tweaks on the Overview/zh-cn/de-de sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

# Drop the existing hyperlinks (Range.Hyperlinks.Delete clears the whole
# sheet's collection in this host) so we can re-add them with the new
# display text, in the same order -> same rId2/rId3 assignment.
$ov.Range("A1").Hyperlinks.Delete()

$ov.Range("A2").Value = "ccbd7795-b258-482a-b40c-ae955711cc8a.md"
$ov.Range("A3").Value = "ffff16413676-89f6-4539-a08f-d0c8e8575bd6.md"

$ov.Hyperlinks.Add($ov.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fd9579c8128bb7c36c118e306cd579c7181ce517/e2e/23f6bd92-b7d1-4908-94ed-2075f84ea54d.md", [Type]::Missing, [Type]::Missing, "e2e\ccbd7795-b258-482a-b40c-ae955711cc8a.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fd9579c8128bb7c36c118e306cd579c7181ce517/e2e/2fe9b83d-e653-41fc-9bec-d0f2cbcbf83b.md", [Type]::Missing, [Type]::Missing, "e2e\ffff16413676-89f6-4539-a08f-d0c8e8575bd6.md") | Out-Null

$ov.Range("E2").Value = "Ready for handoff"
$ov.Range("F2").Value = "Ready for handoff"
$ov.Range("E3").Value = "Ready for handoff"
$ov.Range("F3").Value = "Ready for handoff"

$ov.Range("G2").Value = "2016-08-27 11:02:09"
$ov.Range("G3").Value = "2016-08-27 11:02:09"

$ov.Columns.Item(5).ColumnWidth = 16.33
$ov.Columns.Item(6).ColumnWidth = 16.33

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A1").Hyperlinks.Delete()

$zh.Range("A2").Value = "ccbd7795-b258-482a-b40c-ae955711cc8a.md"
$zh.Range("A3").Value = "ffff16413676-89f6-4539-a08f-d0c8e8575bd6.md"

$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fd9579c8128bb7c36c118e306cd579c7181ce517/e2e/23f6bd92-b7d1-4908-94ed-2075f84ea54d.md", [Type]::Missing, [Type]::Missing, "ccbd7795-b258-482a-b40c-ae955711cc8a.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fd9579c8128bb7c36c118e306cd579c7181ce517/e2e/2fe9b83d-e653-41fc-9bec-d0f2cbcbf83b.md", [Type]::Missing, [Type]::Missing, "ffff16413676-89f6-4539-a08f-d0c8e8575bd6.md") | Out-Null

$zh.Range("C2").Value = "Ready for handoff"
$zh.Range("C3").Value = "Ready for handoff"

$zh.Range("F3").Value = "True"

$zh.Range("G2").Value = "ccbd7795-b258-482a-b40c-ae955711cc8a.b16176464e09a6896749547b167f4b1fdf4f19ca.zh-cn.xlf"
$zh.Range("G3").Value = "ccbd7795-b258-482a-b40c-ae955711cc8a.b16176464e09a6896749547b167f4b1fdf4f19ca.zh-cn.xlf"

$zh.Range("H2").Value = "2016-08-27 11:02:00"
$zh.Range("H3").Value = "2016-08-27 11:02:00"

$zh.Range("I2").Value = ""
$zh.Range("I2").Style = "Normal"
$zh.Range("I3").Value = ""
$zh.Range("I3").Style = "Normal"

$zh.Range("J2").Value = ""
$zh.Range("J3").Value = ""

$zh.Range("K2").Value = "0001-01-01 00:00:00"
$zh.Range("K3").Value = "0001-01-01 00:00:00"

$zh.Columns.Item(3).ColumnWidth = 16.33
$zh.Columns.Item(9).ColumnWidth = 17.83
$zh.Columns.Item(10).ColumnWidth = 20.83

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A1").Hyperlinks.Delete()

$de.Range("A2").Value = "ccbd7795-b258-482a-b40c-ae955711cc8a.md"
$de.Range("A3").Value = "ffff16413676-89f6-4539-a08f-d0c8e8575bd6.md"

$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fd9579c8128bb7c36c118e306cd579c7181ce517/e2e/23f6bd92-b7d1-4908-94ed-2075f84ea54d.md", [Type]::Missing, [Type]::Missing, "ccbd7795-b258-482a-b40c-ae955711cc8a.md") | Out-Null
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fd9579c8128bb7c36c118e306cd579c7181ce517/e2e/2fe9b83d-e653-41fc-9bec-d0f2cbcbf83b.md", [Type]::Missing, [Type]::Missing, "ffff16413676-89f6-4539-a08f-d0c8e8575bd6.md") | Out-Null

$de.Range("C2").Value = "Ready for handoff"
$de.Range("C3").Value = "Ready for handoff"

$de.Range("F3").Value = "True"

$de.Range("G2").Value = "ccbd7795-b258-482a-b40c-ae955711cc8a.b16176464e09a6896749547b167f4b1fdf4f19ca.de-de.xlf"
$de.Range("G3").Value = "ccbd7795-b258-482a-b40c-ae955711cc8a.b16176464e09a6896749547b167f4b1fdf4f19ca.de-de.xlf"

$de.Range("H2").Value = "2016-08-27 11:02:09"
$de.Range("H3").Value = "2016-08-27 11:02:09"

$de.Range("I2").Value = ""
$de.Range("I2").Style = "Normal"
$de.Range("I3").Value = ""
$de.Range("I3").Style = "Normal"

$de.Range("J2").Value = ""
$de.Range("J3").Value = ""

$de.Range("K2").Value = "0001-01-01 00:00:00"
$de.Range("K3").Value = "0001-01-01 00:00:00"

$de.Columns.Item(3).ColumnWidth = 16.33
$de.Columns.Item(9).ColumnWidth = 17.83
$de.Columns.Item(10).ColumnWidth = 20.83
